$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update revised AgTests (H) / AgPosit (I) values for rows 272-329
$ws.Range("H272").Value = 30918
$ws.Range("I272").Value = 1662
$ws.Range("H274").Value = 28351
$ws.Range("I274").Value = 1348
$ws.Range("H275").Value = 28741
$ws.Range("I275").Value = 1239
$ws.Range("H278").Value = 29812
$ws.Range("I278").Value = 2097
$ws.Range("H279").Value = 43382
$ws.Range("H280").Value = 35570
$ws.Range("I280").Value = 2405
$ws.Range("H282").Value = 46812
$ws.Range("I282").Value = 2860
$ws.Range("H285").Value = 40884
$ws.Range("I285").Value = 3432
$ws.Range("H286").Value = 54196
$ws.Range("I286").Value = 4250
$ws.Range("H287").Value = 57649
$ws.Range("I287").Value = 3917
$ws.Range("H288").Value = 56086
$ws.Range("I288").Value = 3959
$ws.Range("H289").Value = 64395
$ws.Range("I289").Value = 3699
$ws.Range("H292").Value = 81241
$ws.Range("I292").Value = 7191
$ws.Range("H293").Value = 81666
$ws.Range("I293").Value = 5794
$ws.Range("H294").Value = 90790
$ws.Range("I294").Value = 5025
$ws.Range("H299").Value = 64055
$ws.Range("I299").Value = 6724
$ws.Range("H300").Value = 70463
$ws.Range("I300").Value = 6917
$ws.Range("H301").Value = 69574
$ws.Range("I301").Value = 5560
$ws.Range("H302").Value = 72555
$ws.Range("I302").Value = 5272
$ws.Range("H306").Value = 70746
$ws.Range("I306").Value = 7172
$ws.Range("H307").Value = 73347
$ws.Range("I307").Value = 6323
$ws.Range("H309").Value = 57315
$ws.Range("I309").Value = 3962
$ws.Range("H310").Value = 90964
$ws.Range("I310").Value = 5190
$ws.Range("H313").Value = 72991
$ws.Range("I313").Value = 3552
$ws.Range("H314").Value = 65084
$ws.Range("I314").Value = 3346
$ws.Range("H315").Value = 66364
$ws.Range("I315").Value = 3012
$ws.Range("H317").Value = 61554
$ws.Range("I317").Value = 2142
$ws.Range("H320").Value = 86293
$ws.Range("I320").Value = 3895
$ws.Range("H321").Value = 90180
$ws.Range("I321").Value = 2796
$ws.Range("H322").Value = 104457
$ws.Range("H323").Value = 150035
$ws.Range("H324").Value = 231685
$ws.Range("I324").Value = 2667
$ws.Range("H325").Value = 673247
$ws.Range("I325").Value = 5517
$ws.Range("H326").Value = 404675
$ws.Range("I326").Value = 3529
$ws.Range("H327").Value = 253652
$ws.Range("I327").Value = 3774
$ws.Range("H328").Value = 189102
$ws.Range("I328").Value = 2751
$ws.Range("H329").Value = 82357
$ws.Range("I329").Value = 1845

# Append new row 330 for 2021-01-28
$ws.Range("A330").Value = 44224
$ws.Range("A330").NumberFormat = "yyyy-mm-dd"
$ws.Range("B330").Value = 246008
$ws.Range("C330").Value = 208406
$ws.Range("D330").Value = 33107
$ws.Range("E330").Value = 11293
$ws.Range("F330").Value = 2581
$ws.Range("G330").Value = 4495
$ws.Range("H330").Value = 66287
$ws.Range("I330").Value = 5098
